# Optimization Updates - Attribute Mapper
#
# The lookup table no longer needs the "reference_codespace_value"
# (always "EPSG") and "reference_version_value" (always "unknown")
# columns, so they are dropped, leaving just original_value /
# reference_code_value. Two additional lookup rows are appended for the
# EPSG:4269 and EPSG:3400 codes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the constant "unknown" column (D / reference_version_value) first,
# then delete the constant "EPSG" column (C / reference_codespace_value)
# so column D's formatting/width shifts left into column C's place -
# mirroring what Excel does when you delete an interior column.
$ws.Range("D1:D10").ClearContents()
$ws.Range("C1").EntireColumn.Delete()

# Append the two new reference rows.
$ws.Range("A11").Value = "EPSG:4269"
$ws.Range("B11").Value = "EPSG: 4269"

$ws.Range("A12").Value = "EPSG:3400"
$ws.Range("B12").Value = "EPSG: 3400"

$ws.Range("B12").Select()
